# Change the table style applied to the table on slide 16
# (the "PLENARY - COMPLETE THE MISSING GAPS" cash-flow table)
# from the default "Table_0" style to the built-in
# "Medium Style 2 - Accent 1" style ({91A1530C-754E-4E84-ACD4-4E96FBA26102}).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

# The table lives in the 3rd shape on the slide (title textbox, picture,
# then the graphic frame holding the table).
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{91A1530C-754E-4E84-ACD4-4E96FBA26102}")
    }
}
